$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usages")

$ws.Range("C5").Value = 99
$ws.Range("C6").Value = 99
$ws.Range("C12").Value = 99
$ws.Range("C13").Value = 99
$ws.Range("C19").Value = 99
$ws.Range("C20").Value = 99
$ws.Range("C26").Value = 99
$ws.Range("C27").Value = 99

foreach ($addr in @("C5","C6","C12","C13","C19","C20","C26","C27")) {
    $ws.Range($addr).Interior.Color = 192
}
